# Auto-generated edit script applying updated probability matrix values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.1918367346938775
$ws.Cells.Item(2, 3).Value = 0.5387755102040817
$ws.Cells.Item(2, 10).Value = 0.0163265306122449
$ws.Cells.Item(2, 16).Value = 0.1551020408163265
$ws.Cells.Item(2, 19).Value = 0.09795918367346938
$ws.Cells.Item(3, 2).Value = 0.007194244604316547
$ws.Cells.Item(3, 3).Value = 0.04316546762589928
$ws.Cells.Item(3, 10).Value = 0.01438848920863309
$ws.Cells.Item(3, 16).Value = 0.7913669064748201
$ws.Cells.Item(3, 19).Value = 0.1438848920863309
$ws.Cells.Item(4, 10).Value = 0.06521739130434782
$ws.Cells.Item(4, 16).Value = 0.6739130434782609
$ws.Cells.Item(4, 19).Value = 0.2608695652173913
$ws.Cells.Item(6, 2).Value = 0.03571428571428571
$ws.Cells.Item(6, 4).Value = 0.00510204081632653
$ws.Cells.Item(6, 5).Value = 0.01020408163265306
$ws.Cells.Item(6, 6).Value = 0.04081632653061224
$ws.Cells.Item(6, 10).Value = 0.25
$ws.Cells.Item(6, 15).Value = 0.02040816326530612
$ws.Cells.Item(6, 17).Value = 0.2091836734693878
$ws.Cells.Item(6, 18).Value = 0.08673469387755102
$ws.Cells.Item(6, 19).Value = 0.3418367346938775
$ws.Cells.Item(7, 2).Value = 0.1
$ws.Cells.Item(7, 4).Value = 0.02777777777777778
$ws.Cells.Item(7, 6).Value = 0.08888888888888889
$ws.Cells.Item(7, 10).Value = 0.08888888888888889
$ws.Cells.Item(7, 15).Value = 0.02777777777777778
$ws.Cells.Item(7, 17).Value = 0.1888888888888889
$ws.Cells.Item(7, 18).Value = 0.06111111111111111
$ws.Cells.Item(7, 19).Value = 0.4166666666666667
$ws.Cells.Item(8, 2).Value = 0.1372549019607843
$ws.Cells.Item(8, 4).Value = 0.0196078431372549
$ws.Cells.Item(8, 6).Value = 0.06162464985994398
$ws.Cells.Item(8, 10).Value = 0.08683473389355742
$ws.Cells.Item(8, 15).Value = 0.01680672268907563
$ws.Cells.Item(8, 17).Value = 0.2240896358543417
$ws.Cells.Item(8, 18).Value = 0.04761904761904762
$ws.Cells.Item(8, 19).Value = 0.4061624649859944
$ws.Cells.Item(9, 2).Value = 0.09937888198757763
$ws.Cells.Item(9, 4).Value = 0.01863354037267081
$ws.Cells.Item(9, 6).Value = 0.06211180124223602
$ws.Cells.Item(9, 10).Value = 0.04968944099378882
$ws.Cells.Item(9, 15).Value = 0.02484472049689441
$ws.Cells.Item(9, 17).Value = 0.1863354037267081
$ws.Cells.Item(9, 18).Value = 0.06211180124223602
$ws.Cells.Item(9, 19).Value = 0.4968944099378882
$ws.Cells.Item(10, 2).Value = 0.09483568075117371
$ws.Cells.Item(10, 4).Value = 0.03098591549295775
$ws.Cells.Item(10, 6).Value = 0.07042253521126761
$ws.Cells.Item(10, 10).Value = 0.1136150234741784
$ws.Cells.Item(10, 15).Value = 0.0215962441314554
$ws.Cells.Item(10, 17).Value = 0.2262910798122066
$ws.Cells.Item(10, 18).Value = 0.05727699530516432
$ws.Cells.Item(10, 19).Value = 0.3849765258215962
$ws.Cells.Item(11, 7).Value = 0.135048231511254
$ws.Cells.Item(11, 10).Value = 0.1028938906752412
$ws.Cells.Item(11, 11).Value = 0.2122186495176849
$ws.Cells.Item(11, 12).Value = 0.5401929260450161
$ws.Cells.Item(11, 19).Value = 0.009646302250803859
$ws.Cells.Item(12, 7).Value = 0.6823529411764706
$ws.Cells.Item(12, 10).Value = 0.2470588235294118
$ws.Cells.Item(12, 12).Value = 0.005882352941176471
$ws.Cells.Item(12, 19).Value = 0.06470588235294118
$ws.Cells.Item(13, 7).Value = 0.5454545454545454
$ws.Cells.Item(13, 10).Value = 0.3409090909090909
$ws.Cells.Item(13, 19).Value = 0.1136363636363636
$ws.Cells.Item(14, 7).Value = 0.6666666666666666
$ws.Cells.Item(14, 19).Value = 0.3333333333333333
$ws.Cells.Item(15, 6).Value = 0.0273972602739726
$ws.Cells.Item(15, 8).Value = 0.1141552511415525
$ws.Cells.Item(15, 9).Value = 0.0867579908675799
$ws.Cells.Item(15, 10).Value = 0.3926940639269406
$ws.Cells.Item(15, 11).Value = 0.0410958904109589
$ws.Cells.Item(15, 13).Value = 0.0136986301369863
$ws.Cells.Item(15, 15).Value = 0.0502283105022831
$ws.Cells.Item(15, 19).Value = 0.273972602739726
$ws.Cells.Item(16, 6).Value = 0.02325581395348837
$ws.Cells.Item(16, 8).Value = 0.1686046511627907
$ws.Cells.Item(16, 9).Value = 0.06395348837209303
$ws.Cells.Item(16, 10).Value = 0.3604651162790697
$ws.Cells.Item(16, 11).Value = 0.1337209302325581
$ws.Cells.Item(16, 13).Value = 0.005813953488372093
$ws.Cells.Item(16, 15).Value = 0.05813953488372093
$ws.Cells.Item(16, 19).Value = 0.186046511627907
$ws.Cells.Item(17, 6).Value = 0.02347417840375587
$ws.Cells.Item(17, 8).Value = 0.1877934272300469
$ws.Cells.Item(17, 9).Value = 0.08215962441314555
$ws.Cells.Item(17, 10).Value = 0.3802816901408451
$ws.Cells.Item(17, 11).Value = 0.1103286384976526
$ws.Cells.Item(17, 13).Value = 0.01643192488262911
$ws.Cells.Item(17, 15).Value = 0.06338028169014084
$ws.Cells.Item(17, 19).Value = 0.136150234741784
$ws.Cells.Item(18, 6).Value = 0.01739130434782609
$ws.Cells.Item(18, 8).Value = 0.1652173913043478
$ws.Cells.Item(18, 9).Value = 0.1043478260869565
$ws.Cells.Item(18, 10).Value = 0.3739130434782609
$ws.Cells.Item(18, 11).Value = 0.1391304347826087
$ws.Cells.Item(18, 13).Value = 0.01739130434782609
$ws.Cells.Item(18, 15).Value = 0.1043478260869565
$ws.Cells.Item(18, 19).Value = 0.0782608695652174
$ws.Cells.Item(19, 6).Value = 0.02079722703639515
$ws.Cells.Item(19, 8).Value = 0.1793760831889082
$ws.Cells.Item(19, 9).Value = 0.07538994800693241
$ws.Cells.Item(19, 10).Value = 0.3500866551126516
$ws.Cells.Item(19, 11).Value = 0.1239168110918544
$ws.Cells.Item(19, 13).Value = 0.0268630849220104
$ws.Cells.Item(19, 14).Value = 0.004332755632582322
$ws.Cells.Item(19, 15).Value = 0.08145580589254767
$ws.Cells.Item(19, 19).Value = 0.1377816291161179
